$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.039.54"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "2.602.59"
$ws.Range("E3").Value = "  -2.44%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'602.69"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("D6").Value = "'145.01"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "'0.584"
$ws.Range("E8").Value = "  -1.00%  "
$ws.Range("D9").Value = "'0.108"
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("D10").Value = "'5.53"
$ws.Range("E10").Value = "  -3.22%  "
$ws.Range("D11").Value = "'0.369"
$ws.Range("E11").Value = "  +3.34%  "
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").Value = "'27.17"
$ws.Range("E13").Value = "  -1.74%  "
$ws.Range("D14").Value = "3.080.69"
$ws.Range("E14").Value = "  -1.96%  "
$ws.Range("D15").Value = "62.920.56"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("D16").Value = "'0.0000145"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "2.623.75"
$ws.Range("E17").Value = "  -1.29%  "
$ws.Range("D18").Value = "'11.42"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("D19").Value = "'4.52"
$ws.Range("E19").Value = "  +2.40%  "
$ws.Range("D20").Value = "'340.74"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "'6.82"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "'5.68"
$ws.Range("E23").Value = "  -1.82%  "
$ws.Range("D24").Value = "'66.37"
$ws.Range("E24").Value = "  -1.72%  "
$ws.Range("D25").Value = "'1.68"
$ws.Range("E25").Value = "  -0.49%  "
$ws.Range("D26").Value = "'9.00"
$ws.Range("E26").Value = "  +5.85%  "
$ws.Range("E27").Value = "  +1.57%  "
$ws.Range("D28").Value = "'549.84"
$ws.Range("E28").Value = "  +1.68%  "
$ws.Range("E29").Value = "  -2.61%  "
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("D31").Value = "'7.80"
$ws.Range("E31").Value = "  -1.41%  "
$ws.Range("D32").Value = "'2.03"
$ws.Range("E32").Value = "  +1.68%  "
$ws.Range("D33").Value = "0.0₃0836"
$ws.Range("E33").Value = "  +2.30%  "
$ws.Range("D34").Value = "'1.74"
$ws.Range("E34").Value = "  -6.36%  "
$ws.Range("D35").Value = "'5.12"
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("D36").Value = "'167.32"
$ws.Range("E36").Value = "  -3.90%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").Value = "'0.400"
$ws.Range("E38").Value = "  -1.16%  "
$ws.Range("D39").Value = "'1.91"
$ws.Range("E39").Value = "  +4.18%  "
$ws.Range("D40").Value = "'18.94"
$ws.Range("E40").Value = "  -1.19%  "
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").Value = "'164.65"
$ws.Range("E42").Value = "  -4.99%  "
$ws.Range("D43").Value = "'3.74"
$ws.Range("E43").Value = "  -1.02%  "
$ws.Range("D44").Value = "'21.63"
$ws.Range("E44").Value = "  -2.85%  "
$ws.Range("D45").Value = "'0.0561"
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("D46").Value = "'0.622"
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("D47").Value = "'0.0243"
$ws.Range("E47").Value = "  +0.77%  "
$ws.Range("D48").Value = "'0.0953"
$ws.Range("E48").Value = "  -1.14%  "
$ws.Range("D49").Value = "'1.89"
$ws.Range("E49").Value = "  +8.78%  "
$ws.Range("D50").Value = "'18.57"
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("D51").Value = "'0.178"
$ws.Range("E51").Value = "  +0.99%  "
